# Add two new departure rows (row 20 and row 21) to the "Main Data" sheet,
# mirroring two more Friday Jan 13 flights (FR5218 to Dublin, FR6112 to Gdansk).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 - flight #19
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(20, 3).Value = "6:40 PM"
$ws.Cells.Item(20, 4).Value = "FR5218"
$ws.Cells.Item(20, 5).Value = "Dublin"
$ws.Cells.Item(20, 6).Value = "(DUB)"
$ws.Cells.Item(20, 7).Value = "Ryanair "
$ws.Cells.Item(20, 8).Value = "B738"
$ws.Cells.Item(20, 9).Value = "(EI-EXE)"
$ws.Cells.Item(20, 10).Value = "7:06 PM"
$ws.Cells.Item(20, 12).Value = "0 hours, 26 minutes"

# Row 21 - flight #20
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(21, 3).Value = "9:20 PM"
$ws.Cells.Item(21, 4).Value = "FR6112"
$ws.Cells.Item(21, 5).Value = "Gdansk"
$ws.Cells.Item(21, 6).Value = "(GDN)"
$ws.Cells.Item(21, 7).Value = "Ryanair "
$ws.Cells.Item(21, 8).Value = "B738"
$ws.Cells.Item(21, 9).Value = "(SP-RSW)"
$ws.Cells.Item(21, 10).Value = "9:45 PM"
$ws.Cells.Item(21, 12).Value = "0 hours, 25 minutes"
